# "Generate Report for handoff" - a new handoff round was produced for b.md:
# its Status flips from "Handed back: in sync with en-US" to "Ready for handoff",
# and new Latest-Handoff-File / Latest-Handoff-Datetime values are recorded
# (with the corresponding hyperlink display text) for both the zh-cn and
# de-de target languages.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# --- Overview sheet: row 3 ("b.md.md") status columns ---
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"

# --- zh-cn sheet: row 3 ("b.md.md") ---
$wsZhCn.Range("B3").Value = "Ready for handoff"
$wsZhCn.Range("C3").Value = "b.md.b3a40d6229ff1a8b48804fcfc66c95922eb78fd0.zh-cn.xlf"
$wsZhCn.Range("D3").Value = "2016-01-14 03:10:35"

# The hyperlink display text for C3 needs to change too, but this host's
# Hyperlink.TextToDisplay / .Address setters (and per-item .Delete()) don't
# mutate the existing entry in place -- they leave stray duplicates behind.
# The only operation that reliably clears hyperlinks is wiping the whole
# collection, so rebuild all of this sheet's hyperlinks from scratch,
# keeping every target exactly as it was except the one display string that
# actually changed.
$wsZhCn.Cells.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/cafddb680a83aa4d7bfb0993a974de43ae9670ea/e2e/a.md.md", "", "", "a.md.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f7103e36fcf5b5328e0c111f6873fabb13cb981a/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/a.md.370104d57010292b5765347db07350cde3a977e8.zh-cn.xlf", "", "", "a.md.370104d57010292b5765347db07350cde3a977e8.zh-cn.xlf") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/b0c73b20a9b43ce4c955dff3d38f3c5a6c003f41/e2e/a.md.md", "", "", "a.md.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/f7613d200ca35b19fbbdd3df3d527e1611c4fa2a/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/a.md.370104d57010292b5765347db07350cde3a977e8.zh-cn.xlf", "", "", "a.md.370104d57010292b5765347db07350cde3a977e8.zh-cn.xlf") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/cafddb680a83aa4d7bfb0993a974de43ae9670ea/e2e/b.md.md", "", "", "b.md.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f7103e36fcf5b5328e0c111f6873fabb13cb981a/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/a.md.370104d57010292b5765347db07350cde3a977e8.zh-cn.xlf", "", "", "b.md.b3a40d6229ff1a8b48804fcfc66c95922eb78fd0.zh-cn.xlf") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/b0c73b20a9b43ce4c955dff3d38f3c5a6c003f41/e2e/a.md.md", "", "", "a.md.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/f7613d200ca35b19fbbdd3df3d527e1611c4fa2a/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/a.md.370104d57010292b5765347db07350cde3a977e8.zh-cn.xlf", "", "", "a.md.370104d57010292b5765347db07350cde3a977e8.zh-cn.xlf") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/cafddb680a83aa4d7bfb0993a974de43ae9670ea/.localization-config", "", "", ".localization-config") | Out-Null

# --- de-de sheet: row 3 ("b.md.md") ---
$wsDeDe.Range("B3").Value = "Ready for handoff"
$wsDeDe.Range("C3").Value = "b.md.b3a40d6229ff1a8b48804fcfc66c95922eb78fd0.de-de.xlf"
$wsDeDe.Range("D3").Value = "2016-01-14 03:10:47"

$wsDeDe.Cells.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/cafddb680a83aa4d7bfb0993a974de43ae9670ea/e2e/a.md.md", "", "", "a.md.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d7a3bbd0db394824eb9be2a98b5e5f32eea1ad36/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/a.md.370104d57010292b5765347db07350cde3a977e8.de-de.xlf", "", "", "a.md.370104d57010292b5765347db07350cde3a977e8.de-de.xlf") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/b61b6cfb72526d059ffcc847a1dafb4924919acc/e2e/a.md.md", "", "", "a.md.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/4e08f728812981257e25ee4396f1d74e576168dc/ol-handback/OpenLocalizationTestOrg/oltest.de-de/yuwzho/a.md.370104d57010292b5765347db07350cde3a977e8.de-de.xlf", "", "", "a.md.370104d57010292b5765347db07350cde3a977e8.de-de.xlf") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/cafddb680a83aa4d7bfb0993a974de43ae9670ea/e2e/b.md.md", "", "", "b.md.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d7a3bbd0db394824eb9be2a98b5e5f32eea1ad36/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/a.md.370104d57010292b5765347db07350cde3a977e8.de-de.xlf", "", "", "b.md.b3a40d6229ff1a8b48804fcfc66c95922eb78fd0.de-de.xlf") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/b61b6cfb72526d059ffcc847a1dafb4924919acc/e2e/a.md.md", "", "", "a.md.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/4e08f728812981257e25ee4396f1d74e576168dc/ol-handback/OpenLocalizationTestOrg/oltest.de-de/yuwzho/a.md.370104d57010292b5765347db07350cde3a977e8.de-de.xlf", "", "", "a.md.370104d57010292b5765347db07350cde3a977e8.de-de.xlf") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/cafddb680a83aa4d7bfb0993a974de43ae9670ea/.localization-config", "", "", ".localization-config") | Out-Null
